# CR 70 x 60 en DTV
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) New merged-cell "version history" column (E2:E5 and E8:E9), each
#    built with the correct top/bottom/middle thin-border formatting.
#    Values are written first (text, in the exact order needed so the
#    sharedStrings table ends up in the same sequence as the target),
#    then formatting is layered on / copied across so cellXfs indices
#    get reused the same way Excel would dedupe them.
# ------------------------------------------------------------------

# -- text values (order matters for sharedStrings uniqueCount order) --
$ws.Range("B8").Value  = "En la lista de Tasa aparecen algunas OT de una misma persona ""duplicadas"" (porque  está separando por modelo. Deben salir juntas"
$ws.Range("B9").Value  = "En DTV agregar botón opción ""Otros Recuperos"" similar a como es en Tasa"
$ws.Range("E8").Value  = "V.5 4/1/2022"
$ws.Range("E2").Value  = "V.2 20/12/2022"
$ws.Range("B7").Value  = "No pueden cerrar sesión para cambiar de usuario"
$ws.Range("E7").Value  = "V.4 3/1/2022"
$ws.Range("E6").Value  = "V.3 29/12/2022"
$ws.Range("B10").Value = "En DTV cuando se guarda pone Código cierre 60 y debe ser 70"

# -- build the "top of merge" border style fresh on E8, then reuse it on E2 --
$ws.Range("E8").Borders.Item(7).LineStyle = 1
$ws.Range("E8").Borders.Item(10).LineStyle = 1
$ws.Range("E8").Borders.Item(8).LineStyle = 1
$ws.Range("E8").Font.Name = "Arial"
$ws.Range("E8").Font.Size = 8
$ws.Range("E8").HorizontalAlignment = -4108
$ws.Range("E8").VerticalAlignment = -4108
$ws.Range("E8").NumberFormat = "mm-dd-yy"

$ws.Range("E8").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# -- build the "bottom of merge" border style fresh on E9, then reuse it on E5 --
$ws.Range("E9").Borders.Item(7).LineStyle = 1
$ws.Range("E9").Borders.Item(10).LineStyle = 1
$ws.Range("E9").Borders.Item(9).LineStyle = 1
$ws.Range("E9").Font.Name = "Arial"
$ws.Range("E9").Font.Size = 8
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").NumberFormat = "mm-dd-yy"

$ws.Range("E9").Copy()
$ws.Range("E5").PasteSpecial(-4122)

# -- build the "middle of merge" border style fresh on E3, then reuse it on E4 --
$ws.Range("E3").Borders.Item(7).LineStyle = 1
$ws.Range("E3").Borders.Item(10).LineStyle = 1
$ws.Range("E3").Font.Name = "Arial"
$ws.Range("E3").Font.Size = 8
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").VerticalAlignment = -4108
$ws.Range("E3").NumberFormat = "mm-dd-yy"

$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# -- merge the two ranges --
$ws.Range("E2:E5").Merge()
$ws.Range("E8:E9").Merge()

# ------------------------------------------------------------------
# 2) New highlighted "Subido a la tienda" dates in column C (style
#    copied from the existing C2 highlighted date cell).
# ------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3) New row data for rows 7-10.
# ------------------------------------------------------------------
$ws.Range("A7").Value = 44564
$ws.Range("C7").Value = 44564

$ws.Range("A8").Value = 44565
$ws.Range("C8").Value = 44565

$ws.Range("A9").Value = 44565
$ws.Range("C9").Value = 44565

$ws.Range("A10").Value = 44565

# Row 8 grows to fit the wrapped two-line request text.
$ws.Rows.Item(8).RowHeight = 20.4

# ------------------------------------------------------------------
# 4) Rows 42 / 43: drop the old special border/format (style 6 / 5)
#    back to the plain style used elsewhere in the sheet (style 3).
# ------------------------------------------------------------------
$ws.Range("A42").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("C43").PasteSpecial(-4122)
$ws.Range("D43").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 5) Two new blank rows (47, 48) matching the style of row 46.
# ------------------------------------------------------------------
$ws.Range("A46:E46").Copy()
$ws.Range("A47:E47").PasteSpecial(-4122)
$ws.Range("A48:E48").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 6) Selection moves to F8.
# ------------------------------------------------------------------
$ws.Range("F8").Select()
